# The document has a single section whose primary (default) and first-page
# headers/footers each carry one inline picture (a logo). This edit renames
# the internal "name" label of those four inline pictures:
#   - the two BTec_Logo-Orange pictures (.jpg):   image1.jpg -> image2.jpg
#   - the two Pearson/PowerPoint logo pictures (.png): image2.png -> image1.png
#
# wdHeaderFooterPrimary   = 1
# wdHeaderFooterFirstPage = 2
#
# Renaming is done by selecting the picture first and then renaming it
# through $word.Selection.InlineShapes - this reliably round-trips the
# name change for pictures that live in both headers and footers.

function Rename-LogoPicture($range, $newName) {
    for ($j = 1; $j -le $range.InlineShapes.Count; $j++) {
        $shp = $range.InlineShapes.Item($j)
        $shp.Select()
        $selShp = $word.Selection.InlineShapes.Item(1)
        $selShp.Name = $newName
    }
}

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Headers (BTec_Logo-Orange, .jpg) ---
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        Rename-LogoPicture $hdr.Range "image2.jpg"
    }
}

# --- Footers (Pearson Edexcel logo, .png) ---
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        Rename-LogoPicture $ftr.Range "image1.png"
    }
}
